# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear forecast cells that no longer produce a value after the bug fix
# (insufficient trailing history for the naive y_0/y_1 forecast at series start)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update recalculated forecast values (precision-level shifts from the corrected calc)
$ws.Range("E3").Value = 7.865470614547321
$ws.Range("E4").Value = 4.26438011980097
$ws.Range("E5").Value = -6.760862998203621
$ws.Range("C6").Value = -0.8792832172735632
$ws.Range("C7").Value = 0.5799958470386724
$ws.Range("C8").Value = 0.9337833426867226
$ws.Range("E8").Value = 2.730731696345212
$ws.Range("C10").Value = 2.791140000794257
$ws.Range("E10").Value = 1.68385714213084
$ws.Range("C11").Value = 0.5930547804883446
$ws.Range("E11").Value = -1.194610791899986
$ws.Range("C12").Value = 0.4451370000809529
$ws.Range("E14").Value = 2.553470871380514
$ws.Range("E15").Value = 7.617133650412167
$ws.Range("E16").Value = -0.4341460075841019
$ws.Range("C17").Value = 1.670328650030162
$ws.Range("E17").Value = 2.037906845818593
$ws.Range("E18").Value = 2.383242923544548
$ws.Range("C19").Value = 2.562791874943349
$ws.Range("C20").Value = 2.2044495746113
$ws.Range("C21").Value = 1.526411006965578
$ws.Range("E21").Value = 0.6601843988560452
$ws.Range("E22").Value = 2.066615940231942
$ws.Range("E23").Value = 0.5447775838346436
$ws.Range("E24").Value = 1.55185774637272
$ws.Range("E25").Value = 3.441981941009331
$ws.Range("E26").Value = 3.086275812215322
$ws.Range("C27").Value = 1.064009474888961
$ws.Range("E27").Value = 0.03490120525226903
$ws.Range("E29").Value = 1.816757311461781
$ws.Range("E30").Value = -1.135072001636317
$ws.Range("C31").Value = 1.600647602405014
$ws.Range("E31").Value = 2.904532120297265
$ws.Range("E32").Value = -5.866344937500012
$ws.Range("C33").Value = -6.098343679991236
$ws.Range("E33").Value = -24.68459749742852
$ws.Range("C34").Value = -3.258619210312896
$ws.Range("E34").Value = 11.44905912635792
$ws.Range("C35").Value = -1.289259938979481
$ws.Range("E35").Value = -1.064625611893855
$ws.Range("C36").Value = -2.616267413525608
$ws.Range("E36").Value = -4.982381489483368
$ws.Range("C38").Value = 0.4255262881966759
$ws.Range("C39").Value = 1.017074315159539
$ws.Range("E39").Value = -2.180963197656882
$ws.Range("C41").Value = 3.634271709196679
$ws.Range("E41").Value = 3.085654629190437
$ws.Range("C43").Value = -0.4898151384455596
$ws.Range("E43").Value = -2.335581849600521
$ws.Range("C44").Value = -1.298607950737285
$ws.Range("C45").Value = -0.1325798828871849
$ws.Range("C46").Value = -0.2814561130375703
$ws.Range("E46").Value = -1.104428907745314
$ws.Range("E47").Value = -1.149192815438882
$ws.Range("C48").Value = -0.994151974263302
$ws.Range("E49").Value = 0.04770910652025506
$ws.Range("C50").Value = -0.6470065423293758
$ws.Range("E50").Value = 3.449881734069282
$ws.Range("C51").Value = 3.540789332106176
$ws.Range("C52").Value = 1.069485063776932
$ws.Range("E52").Value = -2.110726282892139
$ws.Range("C53").Value = 2.107899101287591
